# Generate Report for Handoff
# Replaces the old run GUID/hash tokens with the new ones produced by the
# latest handoff run, and bumps the associated timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "71290707-f9e2-4ffd-9314-b599420e8963"
$newGuid = "d93ad8b6-0682-451b-9c8f-ad09f23c6808"

$oldHash = "c2309a831ac7482eea8ebe1cba6a26582fd5cadd"
$newHash = "fed07a16ff314eefc589939fbd7e4b5d2edd8bde"

# The external hyperlink target itself (stored in the worksheet's .rels)
# is not touched by this handoff run - only the cell text / hyperlink
# display text changes.
$mdHyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f06542ce3130b825ebac517a809f67353d72d51/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# "Overview" sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# A2: plain file name, no hyperlink on this sheet.
$wsOverview.Range("A2").Value = "$newGuid.md"

# B2: path + file name, carries a hyperlink whose display text must be
# refreshed along with the cell text.
$newB2 = "e2e\$newGuid.md"
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $mdHyperlinkAddress, [System.Type]::Missing, [System.Type]::Missing, $newB2) | Out-Null

# G2: "Latest HO Xliff Generate Date" for this handoff run.
$wsOverview.Range("G2").Value = "2016-08-15 18:52:13"

# ---------------------------------------------------------------------
# "zh-cn" sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# A2: file name, carries a hyperlink.
$newA2 = "$newGuid.md"
$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdHyperlinkAddress, [System.Type]::Missing, [System.Type]::Missing, $newA2) | Out-Null

# G2: Latest Handoff File for zh-cn.
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"

# H2: Latest Handoff Datetime for zh-cn.
$wsZhCn.Range("H2").Value = "2016-08-15 18:52:08"

# ---------------------------------------------------------------------
# "de-de" sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# A2: file name, carries a hyperlink.
$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdHyperlinkAddress, [System.Type]::Missing, [System.Type]::Missing, $newA2) | Out-Null

# G2: Latest Handoff File for de-de.
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"

# H2 ("Latest Handoff Datetime" for de-de) is unchanged by this run.
